$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update existing A5/A6 descriptions: append batchsize note ---
$ws.Range("A5").Value = "Just TPR no LSTM in `nphrase embedding layer `nbatchsize = 60"
$ws.Range("A6").Value = "Just LSTM no TPR in `nphrase embedding layer `nbatchsize = 60"

# --- Update F6 pane number (1 -> 3) ---
$ws.Range("F6").Value = 3

# --- Row heights for rows 5 and 6 change from 30 to 45 ---
$ws.Rows.Item(5).RowHeight = 45
$ws.Rows.Item(6).RowHeight = 45

# --- Add new experiment rows 7, 8, 9 ---
# Row 7
$ws.Range("A7").Value = "Just TPR no LSTM in `nphrase embedding layer `nbatchsize = 40"
$ws.Range("B7").Value = "python -m basic.cli --mode train --noload --len_opt --cluster --LSTMandTPR False --justTPR True --batch_size 40 --run_id 2 |& tee /home/hpalangi/QA/TPR_Stuff/Codes/TPR_ver1.0/Log_Files/EXP5.txt"
$ws.Range("C7").Value = "DLT1 / 2"
$ws.Range("D7").Value = "EXP5.txt"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 1
$ws.Range("A7").WrapText = $true
$ws.Rows.Item(7).RowHeight = 45

# Row 8
$ws.Range("A8").Value = "Just LSTM no TPR in `nphrase embedding layer `nbatchsize = 40"
$ws.Range("B8").Value = "python -m basic.cli --mode train --noload --len_opt --cluster --LSTMandTPR False --justLSTM True --batch_size 40 --run_id 3 |& tee /home/hpalangi/QA/TPR_Stuff/Codes/TPR_ver1.0/Log_Files/EXP6.txt"
$ws.Range("C8").Value = "DLT1 / 5"
$ws.Range("D8").Value = "EXP6.txt"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 2
$ws.Range("A8").WrapText = $true
$ws.Rows.Item(8).RowHeight = 45

# Row 9
$ws.Range("A9").Value = "LSTM output concatenated with TPR output in phrase embedding layer. No mixed TPR+LSTM cell for this experiment. `nbatchsize = 40"
$ws.Range("B9").Value = "python -m basic.cli --mode train --noload --len_opt --cluster --LSTMandTPR True --batch_size 40 --run_id 4 |& tee /home/hpalangi/QA/TPR_Stuff/Codes/TPR_ver1.0/Log_Files/EXP7.txt"
$ws.Range("C9").Value = "DLT1 / 7"
$ws.Range("D9").Value = "EXP7.txt"
$ws.Range("E9").Value = 4
$ws.Range("F9").Value = 4
$ws.Range("A9").WrapText = $true
$ws.Rows.Item(9).RowHeight = 105

# --- Update selection to match final state (B9) ---
$ws.Range("B9").Select() | Out-Null

